# Add a new "Montage Pauschale" line item (Zahl / Montage Pauschale / P_Montage)
# as a new row just above the final "Preis / Gesamtpreis / Endpreis" row on the
# "Brix_Zaun_Stab" and "Brix_Schiebe" sheets, and leave the selection/active
# sheet on "Brix_Schiebe" (cell C16), matching the author's edit.

$wb = $excel.ActiveWorkbook

# --- Brix_Zaun_Stab: insert new row 7 (pushes old row 7 -> row 8) ---
$wsZaun = $wb.Worksheets.Item("Brix_Zaun_Stab")
$wsZaun.Rows.Item(7).Insert()
$wsZaun.Cells.Item(7, 1).Value = "Zahl"
$wsZaun.Cells.Item(7, 2).Value = "Montage Pauschale"
$wsZaun.Cells.Item(7, 3).Value = "P_Montage"
$wsZaun.Range("D12").Select()

# --- Brix_Schiebe: insert new row 7 (pushes old row 7 -> row 8) ---
$wsSchiebe = $wb.Worksheets.Item("Brix_Schiebe")
$wsSchiebe.Rows.Item(7).Insert()
$wsSchiebe.Cells.Item(7, 1).Value = "Zahl"
$wsSchiebe.Cells.Item(7, 2).Value = "Montage Pauschale"
$wsSchiebe.Cells.Item(7, 3).Value = "P_Montage"
$wsSchiebe.Range("C16").Select()

# Leave Brix_Schiebe as the active/selected sheet (activeTab moves there).
$wsSchiebe.Activate()
